$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update company names (rows 3 and 4 swapped) ---
$ws.Range('B3').Value = 'Zurich Insurance Group AG (SWX:ZURN)'
$ws.Range('B4').Value = 'Bâloise Holding AG (SWX:BALN)'

# --- Update numeric values ---
# Row 2
$ws.Range('D2').Value = -0.0281
$ws.Range('E2').Value = -0.0524
$ws.Range('F2').Value = 0.04219999999999999
$ws.Range('G2').Value = 0.1359569136741537
$ws.Range('H2').Value = 0.1359569136741537
$ws.Range('I2').Value = 0.08102797600572283
$ws.Range('J2').Value = 0.06734293369177539
$ws.Range('K2').Value = 4162.2
$ws.Range('L2').Value = 0.05544580306629066
$ws.Range('M2').Value = 4073.21
$ws.Range('N2').Value = 0.05221682518072397
$ws.Range('O2').Value = 0.9786194800826487
$ws.Range('P2').Value = 3402.01
$ws.Range('Q2').Value = 0.04361232576593762
$ws.Range('R2').Value = 0.8173586084282352
$ws.Range('S2').Value = 671.1999999999999
$ws.Range('T2').Value = 0.1647840401059606
$ws.Range('U2').Value = 15808.7
$ws.Range('V2').Value = 0.2026608311956691
$ws.Range('W2').Value = 0.07065985181075603
$ws.Range('X2').Value = 0.04786833020029316
$ws.Range('Y2').Value = 0.02279152161046286
$ws.Range('Z2').Value = 1.478960620365189
$ws.Range('AA2').Value = 0.09178949762418726
$ws.Range('AB2').Value = 0.04102086804008107
$ws.Range('AC2').Value = 0.05100534755634102
$ws.Range('AD2').Value = 20340.3
$ws.Range('AE2').Value = 0
$ws.Range('AF2').Value = 20340.3
$ws.Range('AG2').Value = 4531.599999999999
$ws.Range('AH2').Value = 0.2068238667561467
$ws.Range('AI2').Value = 0.2868572250569052
$ws.Range('AJ2').Value = 0.05490366173838979
$ws.Range('AK2').Value = 0.0822451346402002
$ws.Range('AL2').Value = 589.95
$ws.Range('AM2').Value = 589.95
$ws.Range('AN2').Value = 3.083826071136178
$ws.Range('AO2').Value = 10.31036528519366
$ws.Range('AP2').Value = 0.687043269959671
$ws.Range('AQ2').Value = 10.31036528519366

# Row 3
$ws.Range('D3').Value = -0.04389999999999999
$ws.Range('E3').Value = -0.0329
$ws.Range('F3').Value = 0.04219999999999999
$ws.Range('G3').Value = 0.1296313055439917
$ws.Range('H3').Value = 0.1296313055439917
$ws.Range('I3').Value = 0.09272566019294468
$ws.Range('J3').Value = 0.07094095805832135
$ws.Range('K3').Value = 3287
$ws.Range('L3').Value = 0.05949428948940252
$ws.Range('M3').Value = 3352
$ws.Range('N3').Value = 0.05332230405803096
$ws.Range('O3').Value = 1.019774870702769
$ws.Range('P3').Value = 3037
$ws.Range('Q3').Value = 0.04831140734613366
$ws.Range('R3').Value = 0.923942804989352
$ws.Range('S3').Value = 315
$ws.Range('T3').Value = 0.09397374701670644
$ws.Range('U3').Value = 8689
$ws.Range('V3').Value = 0.1382212112053195
$ws.Range('W3').Value = 0.09982082662698533
$ws.Range('X3').Value = 0.04695387245988113
$ws.Range('Y3').Value = 0.0528669541671042
$ws.Range('Z3').Value = 1.491482871257727
$ws.Range('AA3').Value = 0.1058072238145991
$ws.Range('AB3').Value = 0.04096976422770443
$ws.Range('AC3').Value = 0.06483745958689471
$ws.Range('AD3').Value = 15280
$ws.Range('AE3').Value = 0
$ws.Range('AF3').Value = 15280
$ws.Range('AG3').Value = 6591
$ws.Range('AH3').Value = 0.1955389478264208
$ws.Range('AI3').Value = 0.30536181778213
$ws.Range('AJ3').Value = 0.0948973421257235
$ws.Range('AK3').Value = 0.1593954050785973
$ws.Range('AL3').Value = 381
$ws.Range('AM3').Value = 381
$ws.Range('AN3').Value = 2.809339952197095
$ws.Range('AO3').Value = 13.44619422572178
$ws.Range('AP3').Value = 1.211803640375069
$ws.Range('AQ3').Value = 13.44619422572178

# Row 4
$ws.Range('D4').Value = -0.0314
$ws.Range('E4').Value = -0.0479
$ws.Range('G4').Value = 0.2572183931931794
$ws.Range('H4').Value = 0.2572183931931794
$ws.Range('I4').Value = 0.07153165009986262
$ws.Range('J4').Value = 0.05646877887870726
$ws.Range('K4').Value = 503.6
$ws.Range('L4').Value = 0.05813966912571145
$ws.Range('M4').Value = 648.3
$ws.Range('N4').Value = 0.08110846991117227
$ws.Range('O4').Value = 1.287331215250198
$ws.Range('P4').Value = 303.5
$ws.Range('Q4').Value = 0.03797072438383586
$ws.Range('R4').Value = 0.602660841938046
$ws.Range('S4').Value = 344.8
$ws.Range('T4').Value = 0.5318525374055221
$ws.Range('U4').Value = 4128.8
$ws.Range('V4').Value = 0.516551982985112
$ws.Range('W4').Value = 0.07456874213370845
$ws.Range('X4').Value = 0.0487827879407052
$ws.Range('Y4').Value = 0.02578595419300325
$ws.Range('Z4').Value = 1.688709960423449
$ws.Range('AA4').Value = 0.09535938934542224
$ws.Range('AB4').Value = 0.04107197185245769
$ws.Range('AC4').Value = 0.05428741749296454
$ws.Range('AD4').Value = 2493.1
$ws.Range('AE4').Value = 0
$ws.Range('AF4').Value = 2493.1
$ws.Range('AG4').Value = -1635.7
$ws.Range('AH4').Value = 0.2377528347049904
$ws.Range('AI4').Value = 0.2755202404765325
$ws.Range('AJ4').Value = -0.2572947635002282
$ws.Range('AK4').Value = -0.3324661070346959
$ws.Range('AL4').Value = 153.3
$ws.Range('AM4').Value = 153.3
$ws.Range('AN4').Value = 3.550918672553767
$ws.Range('AO4').Value = 4.041748206131768
$ws.Range('AP4').Value = -2.329725110383137
$ws.Range('AQ4').Value = 4.041748206131768

# Row 5
$ws.Range('D5').Value = -0.0248
$ws.Range('E5').Value = -0.05690000000000001
$ws.Range('G5').Value = 0.1472703257998239
$ws.Range('H5').Value = 0.1472703257998239
$ws.Range('I5').Value = 0.1188729087173466
$ws.Range('J5').Value = 0.1075880430075628
$ws.Range('K5').Value = 126.7
$ws.Range('L5').Value = 0.09297035515115938
$ws.Range('M5').Value = 51.9
$ws.Range('N5').Value = 0.03284394380458169
$ws.Range('O5').Value = 0.409629044988161
$ws.Range('P5').Value = 51.9
$ws.Range('Q5').Value = 0.03284394380458169
$ws.Range('R5').Value = 0.409629044988161
$ws.Range('U5').Value = 192
$ws.Range('V5').Value = 0.1215036071383369
$ws.Range('W5').Value = 0.0667509614878036
$ws.Range('X5').Value = 0.04049632828323475
$ws.Range('Y5').Value = 0.02625463320456885
$ws.Range('Z5').Value = 0.8199759326113116
$ws.Range('AA5').Value = 0.08821960590295226
$ws.Range('AB5').Value = 0.04049632828323475
$ws.Range('AC5').Value = 0.04772327761971751
$ws.Range('AG5').Value = -192
$ws.Range('AJ5').Value = -0.1383086010661288
$ws.Range('AK5').Value = -0.1028112449799197
$ws.Range('AL5').Value = 8.65
$ws.Range('AM5').Value = 8.65
$ws.Range('AO5').Value = 18.72832369942196
$ws.Range('AP5').Value = -1.081081081081081
$ws.Range('AQ5').Value = 18.72832369942196

# Row 6
$ws.Range('D6').Value = 0.00221
$ws.Range('E6').Value = -0.0824
$ws.Range('G6').Value = 0.06282289518286331
$ws.Range('H6').Value = 0.06282289518286331
$ws.Range('I6').Value = 0.01817402135958016
$ws.Range('J6').Value = 0.01571824627667969
$ws.Range('K6').Value = 244.9
$ws.Range('L6').Value = 0.02500459455596169
$ws.Range('M6').Value = 21.01
$ws.Range('N6').Value = 0.003772331448065356
$ws.Range('O6').Value = 0.08579011841567986
$ws.Range('P6').Value = 9.609999999999999
$ws.Range('Q6').Value = 0.001725469072627704
$ws.Range('R6').Value = 0.03924050632911392
$ws.Range('S6').Value = 11.4
$ws.Range('T6').Value = 0.5425987624940505
$ws.Range('U6').Value = 2798.9
$ws.Range('V6').Value = 0.5025406230361792
$ws.Range('W6').Value = 0.04300187880809819
$ws.Range('X6').Value = 0.05274200050080584
$ws.Range('Y6').Value = -0.009740121692707653
$ws.Range('Z6').Value = 1.4147539326005
$ws.Range('AA6').Value = 0.02223745073351576
$ws.Range('AB6').Value = 0.04126023328539066
$ws.Range('AC6').Value = -0.0190227825518749
$ws.Range('AD6').Value = 2567.2
$ws.Range('AE6').Value = 0
$ws.Range('AF6').Value = 2567.2
$ws.Range('AG6').Value = -231.7000000000003
$ws.Range('AH6').Value = 0.3155087443312399
$ws.Range('AI6').Value = 0.2630273969795701
$ws.Range('AJ6').Value = -0.04340739630559412
$ws.Range('AK6').Value = -0.03328401304354076
$ws.Range('AL6').Value = 47
$ws.Range('AM6').Value = 47
$ws.Range('AN6').Value = 9.264525442078671
$ws.Range('AO6').Value = 3.787234042553191
$ws.Range('AP6').Value = -0.836160230963552
$ws.Range('AQ6').Value = 3.787234042553191

# --- Remove cell that no longer exists (F4) ---
$ws.Range('F4').ClearContents()